# Insert a new data row at row 71 (this pushes old rows 71..127 down to 72..128,
# matching the target diff which shows every existing row from 71 onward shifted
# by one position and a new row 128 appended at the end).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with its data.
$ws.Cells.Item(71, 1).Value = 3
$ws.Cells.Item(71, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44574
$ws.Cells.Item(71, 5).Value = 5
$ws.Cells.Item(71, 6).Value = 100112026
$ws.Cells.Item(71, 7).Value = "Haba"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 30
$ws.Cells.Item(71, 11).Value = 8000
$ws.Cells.Item(71, 12).Value = 8000
$ws.Cells.Item(71, 13).Value = 8000
$ws.Cells.Item(71, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(71, 16).Value = 320
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
